$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fbn1"
$ws.Cells.Item(2,3).Value = "Itga5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2.0
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 8.970950666666665
$ws.Cells.Item(2,8).Value = 26.912852
$ws.Cells.Item(2,9).Value = 0.02838798528205506
$ws.Cells.Item(2,10).Value = 0.02838798528205506
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 31.22896466666667
$ws.Cells.Item(2,14).Value = 93.686894
$ws.Cells.Item(2,15).Value = 0.2877106972998646
$ws.Cells.Item(2,16).Value = 0.2877106972998646
$ws.Cells.Item(2,17).Value = 280.153501395743
$ws.Cells.Item(2,18).Value = 2521.381512561687
$ws.Cells.Item(2,19).Value = 0.008167527040438357
$ws.Cells.Item(2,20).Value = 0.008167527040438357
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fbn1"
$ws.Cells.Item(3,3).Value = "Itga5"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2.0
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 8.970950666666665
$ws.Cells.Item(3,8).Value = 26.912852
$ws.Cells.Item(3,9).Value = 0.02838798528205506
$ws.Cells.Item(3,10).Value = 0.02838798528205506
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 40.44578266666667
$ws.Cells.Item(3,14).Value = 121.337348
$ws.Cells.Item(3,15).Value = 0.3726247238124506
$ws.Cells.Item(3,16).Value = 0.3726247238124505
$ws.Cells.Item(3,17).Value = 362.8371209773884
$ws.Cells.Item(3,18).Value = 3265.534088796495
$ws.Cells.Item(3,19).Value = 0.01057806517531768
$ws.Cells.Item(3,20).Value = 0.01057806517531768
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fbn1"
$ws.Cells.Item(4,3).Value = "Itga5"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 2.0
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 8.970950666666665
$ws.Cells.Item(4,8).Value = 26.912852
$ws.Cells.Item(4,9).Value = 0.02838798528205506
$ws.Cells.Item(4,10).Value = 0.02838798528205506
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 25.36964133333333
$ws.Cells.Item(4,14).Value = 76.108924
$ws.Cells.Item(4,15).Value = 0.2337290805561598
$ws.Cells.Item(4,16).Value = 0.2337290805561598
$ws.Cells.Item(4,17).Value = 227.5898008323609
$ws.Cells.Item(4,18).Value = 2048.308207491248
$ws.Cells.Item(4,19).Value = 0.006635097698816526
$ws.Cells.Item(4,20).Value = 0.006635097698816526
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Fbn1"
$ws.Cells.Item(5,3).Value = "Itga5"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 2.0
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 8.970950666666665
$ws.Cells.Item(5,8).Value = 26.912852
$ws.Cells.Item(5,9).Value = 0.02838798528205506
$ws.Cells.Item(5,10).Value = 0.02838798528205506
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 11.49855033333333
$ws.Cells.Item(5,14).Value = 34.495651
$ws.Cells.Item(5,15).Value = 0.1059354983315251
$ws.Cells.Item(5,16).Value = 0.1059354983315251
$ws.Cells.Item(5,17).Value = 103.1529277785169
$ws.Cells.Item(5,18).Value = 928.3763500066519
$ws.Cells.Item(5,19).Value = 0.003007295367482505
$ws.Cells.Item(5,20).Value = 0.003007295367482505
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fbn1"
$ws.Cells.Item(6,3).Value = "Itga5"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 272.2141316666667
$ws.Cells.Item(6,8).Value = 816.6423950000001
$ws.Cells.Item(6,9).Value = 0.8614037742994388
$ws.Cells.Item(6,10).Value = 0.8614037742994389
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 31.22896466666667
$ws.Cells.Item(6,14).Value = 93.686894
$ws.Cells.Item(6,15).Value = 0.2877106972998646
$ws.Cells.Item(6,16).Value = 0.2877106972998646
$ws.Cells.Item(6,17).Value = 8500.965499585682
$ws.Cells.Item(6,18).Value = 76508.68949627114
$ws.Cells.Item(6,19).Value = 0.2478350805604267
$ws.Cells.Item(6,20).Value = 0.2478350805604267
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fbn1"
$ws.Cells.Item(7,3).Value = "Itga5"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 272.2141316666667
$ws.Cells.Item(7,8).Value = 816.6423950000001
$ws.Cells.Item(7,9).Value = 0.8614037742994388
$ws.Cells.Item(7,10).Value = 0.8614037742994389
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 40.44578266666667
$ws.Cells.Item(7,14).Value = 121.337348
$ws.Cells.Item(7,15).Value = 0.3726247238124506
$ws.Cells.Item(7,16).Value = 0.3726247238124505
$ws.Cells.Item(7,17).Value = 11009.91360818538
$ws.Cells.Item(7,18).Value = 99089.22247366846
$ws.Cells.Item(7,19).Value = 0.3209803434893309
$ws.Cells.Item(7,20).Value = 0.3209803434893309
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Fbn1"
$ws.Cells.Item(8,3).Value = "Itga5"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 272.2141316666667
$ws.Cells.Item(8,8).Value = 816.6423950000001
$ws.Cells.Item(8,9).Value = 0.8614037742994388
$ws.Cells.Item(8,10).Value = 0.8614037742994389
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 25.36964133333333
$ws.Cells.Item(8,14).Value = 76.108924
$ws.Cells.Item(8,15).Value = 0.2337290805561598
$ws.Cells.Item(8,16).Value = 0.2337290805561598
$ws.Cells.Item(8,17).Value = 6905.974886248109
$ws.Cells.Item(8,18).Value = 62153.77397623299
$ws.Cells.Item(8,19).Value = 0.2013351121546136
$ws.Cells.Item(8,20).Value = 0.2013351121546136
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Fbn1"
$ws.Cells.Item(9,3).Value = "Itga5"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 272.2141316666667
$ws.Cells.Item(9,8).Value = 816.6423950000001
$ws.Cells.Item(9,9).Value = 0.8614037742994388
$ws.Cells.Item(9,10).Value = 0.8614037742994389
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 11.49855033333333
$ws.Cells.Item(9,14).Value = 34.495651
$ws.Cells.Item(9,15).Value = 0.1059354983315251
$ws.Cells.Item(9,16).Value = 0.1059354983315251
$ws.Cells.Item(9,17).Value = 3130.067894413794
$ws.Cells.Item(9,18).Value = 28170.61104972415
$ws.Cells.Item(9,19).Value = 0.09125323809506766
$ws.Cells.Item(9,20).Value = 0.09125323809506768
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Fbn1"
$ws.Cells.Item(10,3).Value = "Itga5"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 0.405826
$ws.Cells.Item(10,8).Value = 1.217478
$ws.Cells.Item(10,9).Value = 0.001284209772536402
$ws.Cells.Item(10,10).Value = 0.001284209772536402
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 31.22896466666667
$ws.Cells.Item(10,14).Value = 93.686894
$ws.Cells.Item(10,15).Value = 0.2877106972998646
$ws.Cells.Item(10,16).Value = 0.2877106972998646
$ws.Cells.Item(10,17).Value = 12.67352581481467
$ws.Cells.Item(10,18).Value = 114.061732333332
$ws.Cells.Item(10,19).Value = 0.0003694808891357486
$ws.Cells.Item(10,20).Value = 0.0003694808891357486
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Fbn1"
$ws.Cells.Item(11,3).Value = "Itga5"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 0.405826
$ws.Cells.Item(11,8).Value = 1.217478
$ws.Cells.Item(11,9).Value = 0.001284209772536402
$ws.Cells.Item(11,10).Value = 0.001284209772536402
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 40.44578266666667
$ws.Cells.Item(11,14).Value = 121.337348
$ws.Cells.Item(11,15).Value = 0.3726247238124506
$ws.Cells.Item(11,16).Value = 0.3726247238124505
$ws.Cells.Item(11,17).Value = 16.41395019648267
$ws.Cells.Item(11,18).Value = 147.725551768344
$ws.Cells.Item(11,19).Value = 0.0004785283118086267
$ws.Cells.Item(11,20).Value = 0.0004785283118086266
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Fbn1"
$ws.Cells.Item(12,3).Value = "Itga5"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 0.405826
$ws.Cells.Item(12,8).Value = 1.217478
$ws.Cells.Item(12,9).Value = 0.001284209772536402
$ws.Cells.Item(12,10).Value = 0.001284209772536402
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 25.36964133333333
$ws.Cells.Item(12,14).Value = 76.108924
$ws.Cells.Item(12,15).Value = 0.2337290805561598
$ws.Cells.Item(12,16).Value = 0.2337290805561598
$ws.Cells.Item(12,17).Value = 10.29566006374133
$ws.Cells.Item(12,18).Value = 92.66094057367201
$ws.Cells.Item(12,19).Value = 0.0003001571693761683
$ws.Cells.Item(12,20).Value = 0.0003001571693761683
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Fbn1"
$ws.Cells.Item(13,3).Value = "Itga5"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 0.405826
$ws.Cells.Item(13,8).Value = 1.217478
$ws.Cells.Item(13,9).Value = 0.001284209772536402
$ws.Cells.Item(13,10).Value = 0.001284209772536402
$ws.Cells.Item(13,11).Value = 3.0
$ws.Cells.Item(13,12).Value = 1.0
$ws.Cells.Item(13,13).Value = 11.49855033333333
$ws.Cells.Item(13,14).Value = 34.495651
$ws.Cells.Item(13,15).Value = 0.1059354983315251
$ws.Cells.Item(13,16).Value = 0.1059354983315251
$ws.Cells.Item(13,17).Value = 4.666410687575334
$ws.Cells.Item(13,18).Value = 41.997696188178
$ws.Cells.Item(13,19).Value = 0.0001360434022158583
$ws.Cells.Item(13,20).Value = 0.0001360434022158583
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Fbn1"
$ws.Cells.Item(14,3).Value = "Itga5"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3.0
$ws.Cells.Item(14,6).Value = 1.0
$ws.Cells.Item(14,7).Value = 34.42132633333333
$ws.Cells.Item(14,8).Value = 103.263979
$ws.Cells.Item(14,9).Value = 0.1089240306459696
$ws.Cells.Item(14,10).Value = 0.1089240306459696
$ws.Cells.Item(14,11).Value = 3.0
$ws.Cells.Item(14,12).Value = 1.0
$ws.Cells.Item(14,13).Value = 31.22896466666667
$ws.Cells.Item(14,14).Value = 93.686894
$ws.Cells.Item(14,15).Value = 0.2877106972998646
$ws.Cells.Item(14,16).Value = 0.2877106972998646
$ws.Cells.Item(14,17).Value = 1074.94238384347
$ws.Cells.Item(14,18).Value = 9674.481454591225
$ws.Cells.Item(14,19).Value = 0.03133860880986373
$ws.Cells.Item(14,20).Value = 0.03133860880986373
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Fbn1"
$ws.Cells.Item(15,3).Value = "Itga5"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3.0
$ws.Cells.Item(15,6).Value = 1.0
$ws.Cells.Item(15,7).Value = 34.42132633333333
$ws.Cells.Item(15,8).Value = 103.263979
$ws.Cells.Item(15,9).Value = 0.1089240306459696
$ws.Cells.Item(15,10).Value = 0.1089240306459696
$ws.Cells.Item(15,11).Value = 3.0
$ws.Cells.Item(15,12).Value = 1.0
$ws.Cells.Item(15,13).Value = 40.44578266666667
$ws.Cells.Item(15,14).Value = 121.337348
$ws.Cells.Item(15,15).Value = 0.3726247238124506
$ws.Cells.Item(15,16).Value = 0.3726247238124505
$ws.Cells.Item(15,17).Value = 1392.19748397641
$ws.Cells.Item(15,18).Value = 12529.77735578769
$ws.Cells.Item(15,19).Value = 0.04058778683599332
$ws.Cells.Item(15,20).Value = 0.04058778683599331
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Fbn1"
$ws.Cells.Item(16,3).Value = "Itga5"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3.0
$ws.Cells.Item(16,6).Value = 1.0
$ws.Cells.Item(16,7).Value = 34.42132633333333
$ws.Cells.Item(16,8).Value = 103.263979
$ws.Cells.Item(16,9).Value = 0.1089240306459696
$ws.Cells.Item(16,10).Value = 0.1089240306459696
$ws.Cells.Item(16,11).Value = 3.0
$ws.Cells.Item(16,12).Value = 1.0
$ws.Cells.Item(16,13).Value = 25.36964133333333
$ws.Cells.Item(16,14).Value = 76.108924
$ws.Cells.Item(16,15).Value = 0.2337290805561598
$ws.Cells.Item(16,16).Value = 0.2337290805561598
$ws.Cells.Item(16,17).Value = 873.2567032942885
$ws.Cells.Item(16,18).Value = 7859.310329648595
$ws.Cells.Item(16,19).Value = 0.02545871353335344
$ws.Cells.Item(16,20).Value = 0.02545871353335344
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Fbn1"
$ws.Cells.Item(17,3).Value = "Itga5"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3.0
$ws.Cells.Item(17,6).Value = 1.0
$ws.Cells.Item(17,7).Value = 34.42132633333333
$ws.Cells.Item(17,8).Value = 103.263979
$ws.Cells.Item(17,9).Value = 0.1089240306459696
$ws.Cells.Item(17,10).Value = 0.1089240306459696
$ws.Cells.Item(17,11).Value = 3.0
$ws.Cells.Item(17,12).Value = 1.0
$ws.Cells.Item(17,13).Value = 11.49855033333333
$ws.Cells.Item(17,14).Value = 34.495651
$ws.Cells.Item(17,15).Value = 0.1059354983315251
$ws.Cells.Item(17,16).Value = 0.1059354983315251
$ws.Cells.Item(17,17).Value = 395.7953533839255
$ws.Cells.Item(17,18).Value = 3562.158180455329
$ws.Cells.Item(17,19).Value = 0.0115389214667591
$ws.Cells.Item(17,20).Value = 0.0115389214667591
